$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '92.792.97'
$ws.Range('E2').Value = '  -5.61%  '
$ws.Range('D3').Value = '3.369.39'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.23'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -8.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '628.24'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.39%  '
$ws.Range('E7').Value = '  -8.63%  '
$ws.Range('E8').Value = '  -10.20%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.933'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -12.66%  '
$ws.Range('D11').Value = '3.366.07'
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('E12').Value = '  -7.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.32'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -12.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.01'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').Value = '92.907.13'
$ws.Range('E15').Value = '  -5.30%  '
$ws.Range('D16').Value = '3.990.54'
$ws.Range('E16').Value = '  -1.33%  '
$ws.Range('E17').Value = '  -6.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.98'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -12.73%  '
$ws.Range('D19').Value = '3.363.23'
$ws.Range('E19').Value = '  -1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -8.81%  '
$ws.Range('E21').Value = '  -3.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '486.02'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.453'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -16.42%  '
$ws.Range('E24').Value = '  -9.32%  '
$ws.Range('E25').Value = '  -8.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.13'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -10.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '89.62'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.52%  '
$ws.Range('D28').Value = '3.538.98'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.45'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.76%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.28'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -8.70%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  -11.02%  '
$ws.Range('E33').Value = '  -9.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('E35').Value = '  -10.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '28.52'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.530'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -7.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '532.02'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.41'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.69%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -7.26%  '
$ws.Range('E42').Value = '  -5.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.876'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.02'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.60'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('E46').Value = '  -6.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.47'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.13'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.96'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.61%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0392'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -8.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.14'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.65%  '
